# Insert a new weekly price record as row 82 on the (single) data sheet.
# This shifts the previous rows 82-100 down to 83-101, growing the used
# range from A1:R100 to A1:R101 - matching the diff's intent of adding one
# new "Camote"/"Paine" observation while preserving every existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 82..100 down by one row.
$ws.Rows.Item(82).Insert()

# Populate the newly-inserted row 82 with the new record's data.
$ws.Cells.Item(82, 1).Value  = 7
$ws.Cells.Item(82, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(82, 3).Value  = "Ñuble"
$ws.Cells.Item(82, 4).Value  = 44511
$ws.Cells.Item(82, 5).Value  = 16
$ws.Cells.Item(82, 6).Value  = 100112045
$ws.Cells.Item(82, 7).Value  = "Zapallo"
$ws.Cells.Item(82, 8).Value  = "Paine"
$ws.Cells.Item(82, 9).Value  = "1a (guarda)"
$ws.Cells.Item(82, 10).Value = 360
$ws.Cells.Item(82, 11).Value = 220
$ws.Cells.Item(82, 12).Value = 250
$ws.Cells.Item(82, 13).Value = 235
$ws.Cells.Item(82, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(82, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(82, 16).Value = 235
$ws.Cells.Item(82, 17).Value = 1
$ws.Cells.Item(82, 18).Value = "Hortaliza"
